$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.2895605203232839
$ws.Range("J2").Value = 0.2895605203232838
$ws.Range("M2").Value = 6.111751666666666
$ws.Range("N2").Value = 18.335255
$ws.Range("O2").Value = 0.6061514841909396
$ws.Range("P2").Value = 0.6061514841909394
$ws.Range("Q2").Value = 8.335608261359443
$ws.Range("R2").Value = 75.02047435223498
$ws.Range("S2").Value = 0.1755175391570593
$ws.Range("T2").Value = 0.1755175391570592
$ws.Range("I3").Value = 0.2895605203232839
$ws.Range("J3").Value = 0.2895605203232838
$ws.Range("O3").Value = 0.2731664420559804
$ws.Range("P3").Value = 0.2731664420559804
$ws.Range("S3").Value = 0.07909821709658986
$ws.Range("T3").Value = 0.07909821709658985
$ws.Range("I4").Value = 0.2895605203232839
$ws.Range("J4").Value = 0.2895605203232838
$ws.Range("M4").Value = 0.568439
$ws.Range("N4").Value = 1.705317
$ws.Range("O4").Value = 0.0563766596410053
$ws.Range("P4").Value = 0.05637665964100529
$ws.Range("Q4").Value = 0.7752744356943334
$ws.Range("R4").Value = 6.977469921249
$ws.Range("S4").Value = 0.01632445489973817
$ws.Range("T4").Value = 0.01632445489973817
$ws.Range("I5").Value = 0.2895605203232839
$ws.Range("J5").Value = 0.2895605203232838
$ws.Range("M5").Value = 0.3689163333333333
$ws.Range("N5").Value = 1.106749
$ws.Range("O5").Value = 0.036588394815171
$ws.Range("P5").Value = 0.036588394815171
$ws.Range("Q5").Value = 0.5031523209058889
$ws.Range("R5").Value = 4.528370888153
$ws.Range("S5").Value = 0.01059455464047466
$ws.Range("T5").Value = 0.01059455464047466
$ws.Range("I6").Value = 0.2895605203232839
$ws.Range("J6").Value = 0.2895605203232838
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2794673333333333
$ws.Range("N6").Value = 0.838402
$ws.Range("O6").Value = 0.02771701929690382
$ws.Range("P6").Value = 0.02771701929690381
$ws.Range("Q6").Value = 0.3811559008882222
$ws.Range("R6").Value = 3.430403107994
$ws.Range("S6").Value = 0.008025754529421969
$ws.Range("T6").Value = 0.008025754529421966
$ws.Range("G7").Value = 1.444951
$ws.Range("H7").Value = 4.334853
$ws.Range("I7").Value = 0.3067756404662893
$ws.Range("J7").Value = 0.3067756404662893
$ws.Range("M7").Value = 6.111751666666666
$ws.Range("N7").Value = 18.335255
$ws.Range("O7").Value = 0.6061514841909396
$ws.Range("P7").Value = 0.6061514841909394
$ws.Range("Q7").Value = 8.831181682501665
$ws.Range("R7").Value = 79.48063514251498
$ws.Range("S7").Value = 0.1859525097822673
$ws.Range("T7").Value = 0.1859525097822673
$ws.Range("G8").Value = 1.444951
$ws.Range("H8").Value = 4.334853
$ws.Range("I8").Value = 0.3067756404662893
$ws.Range("J8").Value = 0.3067756404662893
$ws.Range("O8").Value = 0.2731664420559804
$ws.Range("P8").Value = 0.2731664420559804
$ws.Range("Q8").Value = 3.979834319104
$ws.Range("R8").Value = 35.818508871936
$ws.Range("S8").Value = 0.08380081021562089
$ws.Range("T8").Value = 0.08380081021562089
$ws.Range("G9").Value = 1.444951
$ws.Range("H9").Value = 4.334853
$ws.Range("I9").Value = 0.3067756404662893
$ws.Range("J9").Value = 0.3067756404662893
$ws.Range("M9").Value = 0.568439
$ws.Range("N9").Value = 1.705317
$ws.Range("O9").Value = 0.0563766596410053
$ws.Range("P9").Value = 0.05637665964100529
$ws.Range("Q9").Value = 0.8213665014889999
$ws.Range("R9").Value = 7.392298513400999
$ws.Range("S9").Value = 0.0172949858687194
$ws.Range("T9").Value = 0.0172949858687194
$ws.Range("G10").Value = 1.444951
$ws.Range("H10").Value = 4.334853
$ws.Range("I10").Value = 0.3067756404662893
$ws.Range("J10").Value = 0.3067756404662893
$ws.Range("M10").Value = 0.3689163333333333
$ws.Range("N10").Value = 1.106749
$ws.Range("O10").Value = 0.036588394815171
$ws.Range("P10").Value = 0.036588394815171
$ws.Range("Q10").Value = 0.5330660247663334
$ws.Range("R10").Value = 4.797594222897
$ws.Range("S10").Value = 0.01122442825305754
$ws.Range("T10").Value = 0.01122442825305754
$ws.Range("G11").Value = 1.444951
$ws.Range("H11").Value = 4.334853
$ws.Range("I11").Value = 0.3067756404662893
$ws.Range("J11").Value = 0.3067756404662893
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2794673333333333
$ws.Range("N11").Value = 0.838402
$ws.Range("O11").Value = 0.02771701929690382
$ws.Range("P11").Value = 0.02771701929690381
$ws.Range("Q11").Value = 0.4038166027673333
$ws.Range("R11").Value = 3.634349424906
$ws.Range("S11").Value = 0.008502906346624168
$ws.Range("T11").Value = 0.008502906346624166
$ws.Range("G12").Value = 0.680678
$ws.Range("H12").Value = 2.042034
$ws.Range("I12").Value = 0.1445138481521608
$ws.Range("J12").Value = 0.1445138481521608
$ws.Range("M12").Value = 6.111751666666666
$ws.Range("N12").Value = 18.335255
$ws.Range("O12").Value = 0.6061514841909396
$ws.Range("P12").Value = 0.6061514841909394
$ws.Range("Q12").Value = 4.160134900963333
$ws.Range("R12").Value = 37.44121410867
$ws.Range("S12").Value = 0.08759728354357635
$ws.Range("T12").Value = 0.08759728354357631
$ws.Range("G13").Value = 0.680678
$ws.Range("H13").Value = 2.042034
$ws.Range("I13").Value = 0.1445138481521608
$ws.Range("J13").Value = 0.1445138481521608
$ws.Range("O13").Value = 0.2731664420559804
$ws.Range("P13").Value = 0.2731664420559804
$ws.Range("Q13").Value = 1.874794138112
$ws.Range("R13").Value = 16.873147243008
$ws.Range("S13").Value = 0.03947633372754399
$ws.Range("T13").Value = 0.03947633372754398
$ws.Range("G14").Value = 0.680678
$ws.Range("H14").Value = 2.042034
$ws.Range("I14").Value = 0.1445138481521608
$ws.Range("J14").Value = 0.1445138481521608
$ws.Range("M14").Value = 0.568439
$ws.Range("N14").Value = 1.705317
$ws.Range("O14").Value = 0.0563766596410053
$ws.Range("P14").Value = 0.05637665964100529
$ws.Range("Q14").Value = 0.386923921642
$ws.Range("R14").Value = 3.482315294778
$ws.Range("S14").Value = 0.008147208030686293
$ws.Range("T14").Value = 0.00814720803068629
$ws.Range("G15").Value = 0.680678
$ws.Range("H15").Value = 2.042034
$ws.Range("I15").Value = 0.1445138481521608
$ws.Range("J15").Value = 0.1445138481521608
$ws.Range("M15").Value = 0.3689163333333333
$ws.Range("N15").Value = 1.106749
$ws.Range("O15").Value = 0.036588394815171
$ws.Range("P15").Value = 0.036588394815171
$ws.Range("Q15").Value = 0.2511132319406667
$ws.Range("R15").Value = 2.260019087466
$ws.Range("S15").Value = 0.005287529732450931
$ws.Range("T15").Value = 0.005287529732450929
$ws.Range("G16").Value = 0.680678
$ws.Range("H16").Value = 2.042034
$ws.Range("I16").Value = 0.1445138481521608
$ws.Range("J16").Value = 0.1445138481521608
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2794673333333333
$ws.Range("N16").Value = 0.838402
$ws.Range("O16").Value = 0.02771701929690382
$ws.Range("P16").Value = 0.02771701929690381
$ws.Range("Q16").Value = 0.1902272655186667
$ws.Range("R16").Value = 1.712045389668
$ws.Range("S16").Value = 0.004005493117903269
$ws.Range("T16").Value = 0.004005493117903268
$ws.Range("G17").Value = 1.220628333333333
$ws.Range("H17").Value = 3.661885
$ws.Range("I17").Value = 0.2591499910582661
$ws.Range("J17").Value = 0.2591499910582661
$ws.Range("M17").Value = 6.111751666666666
$ws.Range("N17").Value = 18.335255
$ws.Range("O17").Value = 0.6061514841909396
$ws.Range("P17").Value = 0.6061514841909394
$ws.Range("Q17").Value = 7.460177250630555
$ws.Range("R17").Value = 67.14159525567499
$ws.Range("S17").Value = 0.1570841517080367
$ws.Range("T17").Value = 0.1570841517080367
$ws.Range("G18").Value = 1.220628333333333
$ws.Range("H18").Value = 3.661885
$ws.Range("I18").Value = 0.2591499910582661
$ws.Range("J18").Value = 0.2591499910582661
$ws.Range("O18").Value = 0.2731664420559804
$ws.Range("P18").Value = 0.2731664420559804
$ws.Range("Q18").Value = 3.361981501013334
$ws.Range("R18").Value = 30.25783350912
$ws.Range("S18").Value = 0.07079108101622569
$ws.Range("T18").Value = 0.07079108101622568
$ws.Range("G19").Value = 1.220628333333333
$ws.Range("H19").Value = 3.661885
$ws.Range("I19").Value = 0.2591499910582661
$ws.Range("J19").Value = 0.2591499910582661
$ws.Range("M19").Value = 0.568439
$ws.Range("N19").Value = 1.705317
$ws.Range("O19").Value = 0.0563766596410053
$ws.Range("P19").Value = 0.05637665964100529
$ws.Range("Q19").Value = 0.6938527491716667
$ws.Range("R19").Value = 6.244674742545
$ws.Range("S19").Value = 0.01461001084186144
$ws.Range("T19").Value = 0.01461001084186143
$ws.Range("G20").Value = 1.220628333333333
$ws.Range("H20").Value = 3.661885
$ws.Range("I20").Value = 0.2591499910582661
$ws.Range("J20").Value = 0.2591499910582661
$ws.Range("M20").Value = 0.3689163333333333
$ws.Range("N20").Value = 1.106749
$ws.Range("O20").Value = 0.036588394815171
$ws.Range("P20").Value = 0.036588394815171
$ws.Range("Q20").Value = 0.4503097290961112
$ws.Range("R20").Value = 4.052787561865
$ws.Range("S20").Value = 0.009481882189187876
$ws.Range("T20").Value = 0.009481882189187874
$ws.Range("G21").Value = 1.220628333333333
$ws.Range("H21").Value = 3.661885
$ws.Range("I21").Value = 0.2591499910582661
$ws.Range("J21").Value = 0.2591499910582661
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2794673333333333
$ws.Range("N21").Value = 0.838402
$ws.Range("O21").Value = 0.02771701929690382
$ws.Range("P21").Value = 0.02771701929690381
$ws.Range("Q21").Value = 0.3411257453077778
$ws.Range("R21").Value = 3.07013170777
$ws.Range("S21").Value = 0.007182865302954414
$ws.Range("T21").Value = 0.007182865302954412
